$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 571.41174
$ws.Range("I9").Value = 571.41174
$ws.Range("K9").Value = 571.41174
$ws.Range("M9").Value = -402.41174

# Row 41
$ws.Range("H41").Value = 1277.7273
$ws.Range("J41").Value = 387.5
$ws.Range("L41").Value = 387.5
$ws.Range("N41").Value = -1267.5

# Row 49
$ws.Range("H49").Value = 1019
$ws.Range("J49").Value = 1019
$ws.Range("L49").Value = 3057
$ws.Range("N49").Value = -3329

# Row 58
$ws.Range("H58").Value = 660216.75
$ws.Range("I58").Value = 1515498.5
$ws.Range("J58").Value = 2307.6924
$ws.Range("K58").Value = 4546495.5
$ws.Range("L58").Value = 6923.0772
$ws.Range("M58").Value = -4546345.5
$ws.Range("N58").Value = -7223.0772

# Row 61
$ws.Range("H61").Value = 3992787
$ws.Range("I61").Value = 6944523
$ws.Range("J61").Value = 57139.332
$ws.Range("K61").Value = 20833569
$ws.Range("L61").Value = 171417.996
$ws.Range("M61").Value = -20833397
$ws.Range("N61").Value = -171761.996

# Row 98
$ws.Range("H98").Value = 888.7778
$ws.Range("I98").Value = 599.8570999999999
$ws.Range("J98").Value = 1900
$ws.Range("K98").Value = 599.8570999999999
$ws.Range("L98").Value = 1900
$ws.Range("M98").Value = 898.1429000000001
$ws.Range("N98").Value = -4896

# Row 121
$ws.Range("H121").Value = 1060.25
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1060.25
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 3180.75
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -6674.75

# Row 122
$ws.Range("H122").Value = 888.7778
$ws.Range("I122").Value = 599.8570999999999
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 1799.5713
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = 650.4287000000002
$ws.Range("N122").Value = -10600

$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 1256249.9
$ws.Range("I6").Value = 3337666.8
$ws.Range("K6").Value = 3337666.8
$ws.Range("M6").Value = -3337493.8

# Row 32
$ws.Range("H32").Value = 29171.258
$ws.Range("I32").Value = 4900.0815
$ws.Range("J32").Value = 120654.92
$ws.Range("K32").Value = 4900.0815
$ws.Range("L32").Value = 120654.92
$ws.Range("M32").Value = -4613.0815
$ws.Range("N32").Value = -121228.92

$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 1000
$ws.Range("I17").Value = 1000
$ws.Range("K17").Value = 1000
$ws.Range("M17").Value = -826

# Row 41
$ws.Range("H41").Value = 10331.875
$ws.Range("J41").Value = 12505
$ws.Range("L41").Value = 12505
$ws.Range("N41").Value = -13361

# Row 50
$ws.Range("H50").Value = 13420
$ws.Range("J50").Value = 13420
$ws.Range("L50").Value = 13420
$ws.Range("N50").Value = -14670

# Row 51
$ws.Range("H51").Value = 7927.8
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

# Row 59
$ws.Range("H59").Value = 21388
$ws.Range("J59").Value = 21388
$ws.Range("L59").Value = 21388
$ws.Range("N59").Value = -23678

# Row 61
$ws.Range("H61").Value = 7927.8
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

# Row 68
$ws.Range("H68").Value = 14356.357
$ws.Range("J68").Value = 14356.357
$ws.Range("L68").Value = 14356.357
$ws.Range("N68").Value = -15854.357

# Row 71
$ws.Range("H71").Value = 14356.357
$ws.Range("J71").Value = 14356.357
$ws.Range("L71").Value = 43069.071
$ws.Range("N71").Value = -50557.071

# Row 74
$ws.Range("H74").Value = 22541.385
$ws.Range("J74").Value = 22541.385
$ws.Range("L74").Value = 22541.385
$ws.Range("N74").Value = -24289.385

# Row 75
$ws.Range("H75").Value = 30999
$ws.Range("J75").Value = 30999
$ws.Range("L75").Value = 30999
$ws.Range("N75").Value = -32995

# Row 77
$ws.Range("H77").Value = 22541.385
$ws.Range("J77").Value = 22541.385
$ws.Range("L77").Value = 67624.155
$ws.Range("N77").Value = -76360.155

# Row 78
$ws.Range("H78").Value = 30999
$ws.Range("J78").Value = 30999
$ws.Range("L78").Value = 92997
$ws.Range("N78").Value = -102981

# Row 122
$ws.Range("H122").Value = 3394.8276
$ws.Range("I122").Value = 3633.25
$ws.Range("J122").Value = 3101.3845
$ws.Range("K122").Value = 10899.75
$ws.Range("L122").Value = 9304.1535
$ws.Range("M122").Value = -8449.75
$ws.Range("N122").Value = -14204.1535

$ws = $wb.Worksheets.Item("CUL")
# Row 69
$ws.Range("H69").Value = 2147.7144
$ws.Range("J69").Value = 2300
$ws.Range("L69").Value = 6900
$ws.Range("N69").Value = -8522

# Row 72
$ws.Range("H72").Value = 2147.7144
$ws.Range("J72").Value = 2300
$ws.Range("L72").Value = 20700
$ws.Range("N72").Value = -28812

# Row 131
$ws.Range("H131").Value = 845.8182
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 845.8182
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2537.4546
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -12617.4546

# Row 132
$ws.Range("H132").Value = 2600
$ws.Range("I132").Value = 2533.3333
$ws.Range("K132").Value = 22799.9997
$ws.Range("M132").Value = -20269.9997

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 1264576.2
$ws.Range("I5").Value = 5000000
$ws.Range("J5").Value = 19435
$ws.Range("K5").Value = 5000000
$ws.Range("L5").Value = 19435
$ws.Range("M5").Value = -4999888
$ws.Range("N5").Value = -19659

# Row 80
$ws.Range("H80").Value = 3077.2
$ws.Range("I80").Value = 4980
$ws.Range("J80").Value = 2601.5
$ws.Range("K80").Value = 4980
$ws.Range("L80").Value = 2601.5
$ws.Range("M80").Value = -3982
$ws.Range("N80").Value = -4597.5

# Row 83
$ws.Range("H83").Value = 3077.2
$ws.Range("I83").Value = 4980
$ws.Range("J83").Value = 2601.5
$ws.Range("K83").Value = 24900
$ws.Range("L83").Value = 13007.5
$ws.Range("M83").Value = -19908
$ws.Range("N83").Value = -22991.5

# Row 102
$ws.Range("H102").Value = 302528.44
$ws.Range("I102").Value = 3098.4614
$ws.Range("J102").Value = 858612.7
$ws.Range("K102").Value = 3098.4614
$ws.Range("L102").Value = 858612.7
$ws.Range("M102").Value = -1476.4614
$ws.Range("N102").Value = -861856.7

# Row 122
$ws.Range("H122").Value = 3510.52
$ws.Range("I122").Value = 2893.15
$ws.Range("K122").Value = 8679.450000000001
$ws.Range("M122").Value = -6229.450000000001

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 346538.47
$ws.Range("I2").Value = 346538.47
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 346538.47
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -346426.47
$ws.Range("N2").ClearContents()

# Row 40
$ws.Range("H40").Value = 102200
$ws.Range("I40").Value = 334833.34
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 334833.34
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -334697.34
$ws.Range("N40").Value = -2772

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 24350
$ws.Range("I2").Value = 17250
$ws.Range("K2").Value = 17250
$ws.Range("M2").Value = -17138

# Row 122
$ws.Range("H122").Value = 1241.4348
$ws.Range("I122").Value = 1005.8461
$ws.Range("J122").Value = 1547.7
$ws.Range("K122").Value = 3017.5383
$ws.Range("L122").Value = 4643.1
$ws.Range("M122").Value = -567.5383000000002
$ws.Range("N122").Value = -9543.1

# Row 126
$ws.Range("H126").Value = 2388.5715
$ws.Range("I126").Value = 2544
$ws.Range("K126").Value = 7632
$ws.Range("M126").Value = -5162
